# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (header row 1, G1 = "K").
# Update rows 2-10 with the recalculated K values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 2
